$wb = $excel.ActiveWorkbook

# Sheet 1: "Weekly Quantity"
# Delete row 2 entirely so rows 3-5 shift up to 2-4, matching the new data
# (45130.99999999999/1, 45144.99999999999/2, 45158.99999999999/1).
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(2).Delete()

# Sheet 2: "Monthly Trend"
# Row 2's quantity changes from 2 to 1.
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B2").Value = 1
